# Add ROIs matching step to parcellation stage.
#
# Inserts two new parameter rows after "parcellation.templateScript"
# (row 31):
#   - parcellation.matchROIs
#   - parcellation.lutFile (replacing the old collect_region_properties.lutFile
#     row, which is removed further down the sheet).
#
# The rest of the sheet (rows 1-31 and the old rows 32-34/36-52, which are
# simply pushed down) remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two fresh rows right after row 31 (parcellation.templateScript) ---
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(33).Insert()

# Row insertion at this location pulls in the column-level "text" style
# (style index 1) that is defined for columns F and G, which leaves a stray,
# content-less styled cell behind in F32 (row 32 has no F value). Clear it so
# the saved sheet doesn't carry an empty/unused cell.
$ws.Range("F32").Clear()

# --- Row 32: parcellation.matchROIs ---
$ws.Range("A32").Value = "parcellation.matchROIs"
$ws.Range("D32").Value = "parcellation"
$ws.Range("E32").Value = "logical"
$ws.Range("G32").Value = "standard"
$ws.Range("H32").Value = "Flag whether the parcellation step should reassign the ROIs in the parcellationFile to match the template's color lookup table."

# --- Row 33: parcellation.lutFile ---
$ws.Range("A33").Value = "parcellation.lutFile"
$ws.Range("B33").Value = "parcellation,collect_region_properties"
$ws.Range("E33").Value = "char "
$ws.Range("F33").Value = "isfile nonempty"
$ws.Range("G33").Value = "standard"
$ws.Range("H33").Value = "Freesurfer's color lookup table of the template"

# --- Remove the now-redundant collect_region_properties.lutFile row ---
# After inserting the two rows above, that row (originally row 35) has moved
# down to row 37.
$ws.Rows.Item(37).Delete()

# Restore the selection/view state seen in the final workbook.
$ws.Range("B33").Select()
